$d = $word.ActiveDocument

# 1. Update the footer timestamp (section 1's footer part is shared across
#    sections in this document, so updating it once is sufficient).
$ftr = $d.Sections.Item(1).Footers.Item(1)
$ftr.Range.Find.Execute(
    "2025-06-30 12:12Z / ", $true, $false, $false, $false, $false,
    $true, 1, $false, "2025-07-02 02:48Z / ", 2) | Out-Null

# 2. Add the new character styles: b, i, sub, sup, u
$wdStyleTypeCharacter = 2

$sB = $d.Styles.Add("b", $wdStyleTypeCharacter)
$sB.BaseStyle = "DefaultParagraphFont"
$sB.Priority = 1
$sB.QuickStyle = $true
$sB.Font.Bold = $true

$sI = $d.Styles.Add("i", $wdStyleTypeCharacter)
$sI.BaseStyle = "DefaultParagraphFont"
$sI.Priority = 1
$sI.QuickStyle = $true
$sI.Font.Italic = $true

$sSub = $d.Styles.Add("sub", $wdStyleTypeCharacter)
$sSub.BaseStyle = "DefaultParagraphFont"
$sSub.Priority = 1
$sSub.QuickStyle = $true
$sSub.Font.Subscript = $true

$sSup = $d.Styles.Add("sup", $wdStyleTypeCharacter)
$sSup.BaseStyle = "DefaultParagraphFont"
$sSup.Priority = 1
$sSup.QuickStyle = $true
$sSup.Font.Superscript = $true

$sU = $d.Styles.Add("u", $wdStyleTypeCharacter)
$sU.BaseStyle = "DefaultParagraphFont"
$sU.Priority = 1
$sU.QuickStyle = $true
$sU.Font.Underline = 1
